# Auto-generated Excel COM-interop script applying the "Spriggan_Profits" leve-profit
# data refresh described in the commit diff (scheduled price-data runner).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on affected rows
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1200
$ws.Range("I4").Value = 1250.5
$ws.Range("J4").Value = 1099
$ws.Range("K4").Value = 1250.5
$ws.Range("L4").Value = 1099
$ws.Range("M4").Value = -1136.5
$ws.Range("N4").Value = -1327
$ws.Range("H17").Value = 119290.47
$ws.Range("J17").Value = 119290.47
$ws.Range("L17").Value = 357871.41
$ws.Range("N17").Value = -358207.41
$ws.Range("H18").Value = 3000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 3000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -3568
$ws.Range("H40").Value = 3706947
$ws.Range("I40").Value = 2088.1667
$ws.Range("K40").Value = 2088.1667
$ws.Range("M40").Value = -1913.1667
$ws.Range("H43").Value = 4947.375
$ws.Range("I43").Value = 4880
$ws.Range("J43").Value = 5059.6665
$ws.Range("K43").Value = 4880
$ws.Range("L43").Value = 5059.6665
$ws.Range("M43").Value = -4811
$ws.Range("N43").Value = -5197.6665
$ws.Range("H53").Value = 334.6842
$ws.Range("I53").Value = 263.44446
$ws.Range("K53").Value = 263.44446
$ws.Range("M53").Value = 373.55554
$ws.Range("H80").Value = 602.4286
$ws.Range("I80").Value = 727.3333
$ws.Range("K80").Value = 2181.9999
$ws.Range("M80").Value = -1183.9999
$ws.Range("H83").Value = 602.4286
$ws.Range("I83").Value = 727.3333
$ws.Range("K83").Value = 6545.9997
$ws.Range("M83").Value = -1553.9997
$ws.Range("H116").Value = 5375.4
$ws.Range("I116").Value = 5077.5386
$ws.Range("K116").Value = 5077.5386
$ws.Range("M116").Value = -1635.5386
$ws.Range("H129").Value = 1467.7333
$ws.Range("I129").Value = 1158.8889
$ws.Range("J129").Value = 1931
$ws.Range("K129").Value = 3476.6667
$ws.Range("L129").Value = 5793
$ws.Range("M129").Value = 1523.3333
$ws.Range("N129").Value = -15793
$ws.Range("H132").Value = 2903.9565
$ws.Range("I132").Value = 2847.1904
$ws.Range("K132").Value = 8541.5712
$ws.Range("M132").Value = -6011.5712
$ws.Range("H137").Value = 1687.2
$ws.Range("I137").Value = 1127.0714
$ws.Range("K137").Value = 3381.2142
$ws.Range("M137").Value = -831.2142000000003

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 45459470
$ws.Range("I74").Value = 45459470
$ws.Range("K74").Value = 45459470
$ws.Range("M74").Value = -45458596
$ws.Range("H77").Value = 45459470
$ws.Range("I77").Value = 45459470
$ws.Range("K77").Value = 227297350
$ws.Range("M77").Value = -227292982

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 9500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 9500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 9500
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -9970
$ws.Range("H86").Value = 3359.3635
$ws.Range("I86").Value = 3636.2856
$ws.Range("J86").Value = 2874.75
$ws.Range("K86").Value = 3636.2856
$ws.Range("L86").Value = 2874.75
$ws.Range("M86").Value = -2513.2856
$ws.Range("N86").Value = -5120.75
$ws.Range("H89").Value = 3359.3635
$ws.Range("I89").Value = 3636.2856
$ws.Range("J89").Value = 2874.75
$ws.Range("K89").Value = 18181.428
$ws.Range("L89").Value = 14373.75
$ws.Range("M89").Value = -12565.428
$ws.Range("N89").Value = -25605.75
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H105").Value = 1947.3889
$ws.Range("I105").Value = 1941.625
$ws.Range("K105").Value = 1941.625
$ws.Range("M105").Value = -194.625
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 279.66666
$ws.Range("I7").Value = 21.333334
$ws.Range("J7").Value = 408.83334
$ws.Range("K7").Value = 21.333334
$ws.Range("L7").Value = 408.83334
$ws.Range("M7").Value = 91.66666599999999
$ws.Range("N7").Value = -634.83334
$ws.Range("H16").Value = 1210398.5
$ws.Range("I16").Value = 1812431.1
$ws.Range("J16").Value = 6333
$ws.Range("K16").Value = 1812431.1
$ws.Range("L16").Value = 6333
$ws.Range("M16").Value = -1812144.1
$ws.Range("N16").Value = -6907
$ws.Range("H58").Value = 33341382
$ws.Range("I58").Value = 38470516
$ws.Range("K58").Value = 38470516
$ws.Range("M58").Value = -38470313
$ws.Range("H103").Value = 29713.715
$ws.Range("I103").Value = 22999.5
$ws.Range("K103").Value = 22999.5
$ws.Range("M103").Value = -21827.5
$ws.Range("H113").Value = 1210398.5
$ws.Range("I113").Value = 1812431.1
$ws.Range("J113").Value = 6333
$ws.Range("K113").Value = 1812431.1
$ws.Range("L113").Value = 6333
$ws.Range("M113").Value = -1810261.1
$ws.Range("N113").Value = -10673
$ws.Range("H134").Value = 15627123
$ws.Range("J134").Value = 2998.3333
$ws.Range("L134").Value = 8994.999899999999
$ws.Range("N134").Value = -14064.9999
$ws.Range("H136").Value = 33341382
$ws.Range("I136").Value = 38470516
$ws.Range("K136").Value = 115411548
$ws.Range("M136").Value = -115408998

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 174.33333
$ws.Range("J38").Value = 134.33333
$ws.Range("L38").Value = 402.99999
$ws.Range("N38").Value = -1096.99999
$ws.Range("H109").Value = 1871.8572
$ws.Range("I109").Value = 1871.8572
$ws.Range("K109").Value = 5615.571599999999
$ws.Range("M109").Value = -4575.571599999999
$ws.Range("H112").Value = 12951.111
$ws.Range("I112").Value = 4140.25
$ws.Range("K112").Value = 12420.75
$ws.Range("M112").Value = -11312.75
$ws.Range("H113").Value = 63807.812
$ws.Range("J113").Value = 972.625
$ws.Range("L113").Value = 2917.875
$ws.Range("N113").Value = -7257.875

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2466.3333
$ws.Range("I16").Value = 1224.6666
$ws.Range("K16").Value = 1224.6666
$ws.Range("M16").Value = -1054.6666
$ws.Range("H22").Value = 3599.9092
$ws.Range("I22").Value = 3622.111
$ws.Range("K22").Value = 3622.111
$ws.Range("M22").Value = -3327.111
$ws.Range("H27").Value = 3599.9092
$ws.Range("I27").Value = 3622.111
$ws.Range("K27").Value = 3622.111
$ws.Range("M27").Value = -3515.111
$ws.Range("H40").Value = 6581.143
$ws.Range("I40").Value = 6533
$ws.Range("K40").Value = 6533
$ws.Range("M40").Value = -6397
$ws.Range("H46").Value = 1167
$ws.Range("I46").Value = 1167
$ws.Range("K46").Value = 1167
$ws.Range("M46").Value = -979
$ws.Range("H68").Value = 2273863.5
$ws.Range("I68").Value = 2273863.5
$ws.Range("K68").Value = 2273863.5
$ws.Range("M68").Value = -2273114.5
$ws.Range("H71").Value = 2273863.5
$ws.Range("I71").Value = 2273863.5
$ws.Range("K71").Value = 11369317.5
$ws.Range("M71").Value = -11365573.5
$ws.Range("H82").Value = 2999.8333
$ws.Range("I82").Value = 2750
$ws.Range("K82").Value = 2750
$ws.Range("M82").Value = -2389
$ws.Range("H85").Value = 2999.8333
$ws.Range("I85").Value = 2750
$ws.Range("K85").Value = 2750
$ws.Range("M85").Value = -1502
$ws.Range("H122").Value = 3565.138
$ws.Range("I122").Value = 3565.138
$ws.Range("K122").Value = 10695.414
$ws.Range("M122").Value = -8245.414000000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 80344
$ws.Range("J106").Value = 80344
$ws.Range("L106").Value = 80344
$ws.Range("N106").Value = -82868
$ws.Range("H136").Value = 10206530
$ws.Range("I136").Value = 10871959
$ws.Range("K136").Value = 32615877
$ws.Range("M136").Value = -32613327

Write-Host "Applied all Spriggan_Profits updates"
